$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 181.12
$ws.Range("I15").Value = 181.12
$ws.Range("K15").Value = 543.36
$ws.Range("M15").Value = -374.36
$ws.Range("H19").Value = 8112.5
$ws.Range("J19").Value = 13000.25
$ws.Range("L19").Value = 13000.25
$ws.Range("N19").Value = -13350.25
$ws.Range("H33").Value = 634.6667
$ws.Range("I33").Value = 713.1667
$ws.Range("J33").Value = 320.66666
$ws.Range("K33").Value = 713.1667
$ws.Range("L33").Value = 320.66666
$ws.Range("M33").Value = -484.1667
$ws.Range("N33").Value = -778.66666
$ws.Range("H40").Value = 3700
$ws.Range("I40").Value = 1350
$ws.Range("K40").Value = 1350
$ws.Range("M40").Value = -1175
$ws.Range("H61").Value = 1422.2667
$ws.Range("I61").Value = 485.63635
$ws.Range("J61").Value = 3998
$ws.Range("K61").Value = 1456.90905
$ws.Range("L61").Value = 11994
$ws.Range("M61").Value = -1284.90905
$ws.Range("N61").Value = -12338
$ws.Range("H64").Value = 2968.923
$ws.Range("I64").Value = 2866.6667
$ws.Range("K64").Value = 2866.6667
$ws.Range("M64").Value = -2618.6667
$ws.Range("H67").Value = 2968.923
$ws.Range("I67").Value = 2866.6667
$ws.Range("K67").Value = 2866.6667
$ws.Range("M67").Value = -2008.6667
$ws.Range("H74").Value = 2950
$ws.Range("J74").Value = 3266.6667
$ws.Range("L74").Value = 3266.6667
$ws.Range("N74").Value = -5138.6667
$ws.Range("H77").Value = 2950
$ws.Range("J77").Value = 3266.6667
$ws.Range("L77").Value = 16333.3335
$ws.Range("N77").Value = -25693.3335
$ws.Range("H98").Value = 1719.8
$ws.Range("I98").Value = 1585.625
$ws.Range("J98").Value = 2256.5
$ws.Range("K98").Value = 1585.625
$ws.Range("L98").Value = 2256.5
$ws.Range("M98").Value = -87.625
$ws.Range("N98").Value = -5252.5
$ws.Range("H100").Value = 2979.8
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2979.8
$ws.Range("K100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("M100").Value = 2979.8
$ws.Range("N100").Value = -4061.8
$ws.Range("H122").Value = 1719.8
$ws.Range("I122").Value = 1585.625
$ws.Range("J122").Value = 2256.5
$ws.Range("K122").Value = 4756.875
$ws.Range("L122").Value = 6769.5
$ws.Range("M122").Value = -2306.875
$ws.Range("N122").Value = -11669.5
$ws.Range("H126").Value = 61150
$ws.Range("J126").Value = 61150
$ws.Range("L126").Value = 61150
$ws.Range("N126").Value = -71030
$ws.Range("H132").Value = 4135.2915
$ws.Range("I132").Value = 4184.6523
$ws.Range("K132").Value = 12553.9569
$ws.Range("M132").Value = -10023.9569
$ws.Range("H140").Value = 76823.69
$ws.Range("J140").Value = 76823.69
$ws.Range("L140").Value = 76823.69
$ws.Range("N140").Value = -87183.69

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1849.8667
$ws.Range("I45").Value = 978.2857
$ws.Range("J45").Value = 2612.5
$ws.Range("K45").Value = 978.2857
$ws.Range("L45").Value = 2612.5
$ws.Range("M45").Value = -601.2857
$ws.Range("N45").Value = -3366.5
$ws.Range("H132").Value = 1376399.8
$ws.Range("I132").Value = 2505.2927
$ws.Range("K132").Value = 7515.8781
$ws.Range("M132").Value = -4985.8781

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H125").Value = 42000
$ws.Range("J125").Value = 42000
$ws.Range("L125").Value = 42000
$ws.Range("N125").Value = -51840
$ws.Range("H134").Value = 2831.913
$ws.Range("I134").Value = 2776.1282
$ws.Range("J134").Value = 3142.7144
$ws.Range("K134").Value = 8328.384600000001
$ws.Range("L134").Value = 9428.143199999999
$ws.Range("M134").Value = -5793.384600000001
$ws.Range("N134").Value = -14498.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2400.6
$ws.Range("I99").Value = 1986
$ws.Range("J99").Value = 2531.5264
$ws.Range("K99").Value = 1986
$ws.Range("L99").Value = 2531.5264
$ws.Range("M99").Value = -488
$ws.Range("N99").Value = -5527.526400000001
$ws.Range("H106").Value = 50447.332
$ws.Range("J106").Value = 50447.332
$ws.Range("L106").Value = 50447.332
$ws.Range("N106").Value = -52971.332
$ws.Range("H118").Value = 78513.96000000001
$ws.Range("J118").Value = 78513.96000000001
$ws.Range("L118").Value = 78513.96000000001
$ws.Range("N118").Value = -81827.96000000001
$ws.Range("H126").Value = 2400.6
$ws.Range("I126").Value = 1986
$ws.Range("J126").Value = 2531.5264
$ws.Range("K126").Value = 5958
$ws.Range("L126").Value = 7594.5792
$ws.Range("M126").Value = -3488
$ws.Range("N126").Value = -12534.5792
$ws.Range("H132").Value = 15505888
$ws.Range("I132").Value = 20835092
$ws.Range("J132").Value = 8774262
$ws.Range("K132").Value = 62505276
$ws.Range("L132").Value = 26322786
$ws.Range("M132").Value = -62502746
$ws.Range("N132").Value = -26327846
$ws.Range("H134").Value = 10423522
$ws.Range("I134").Value = 11911358
$ws.Range("K134").Value = 35734074
$ws.Range("M134").Value = -35731539
$ws.Range("H141").Value = 84474.30499999999
$ws.Range("J141").Value = 81173.57000000001
$ws.Range("L141").Value = 81173.57000000001
$ws.Range("N141").Value = -91533.57000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1167121.9
$ws.Range("J129").Value = 1517091.2
$ws.Range("L129").Value = 4551273.6
$ws.Range("N129").Value = -4561273.6
$ws.Range("H139").Value = 374436.25
$ws.Range("I139").Value = 528210.1
$ws.Range("J139").Value = 9223.375
$ws.Range("K139").Value = 1584630.3
$ws.Range("L139").Value = 27670.125
$ws.Range("M139").Value = -1579490.3
$ws.Range("N139").Value = -37950.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 2003000
$ws.Range("I3").Value = 5000000
$ws.Range("J3").Value = 504500
$ws.Range("K3").Value = 5000000
$ws.Range("L3").Value = 504500
$ws.Range("M3").Value = -4999884
$ws.Range("N3").Value = -504732
$ws.Range("H14").Value = 12222550
$ws.Range("I14").Value = 12222550
$ws.Range("K14").Value = 12222550
$ws.Range("M14").Value = -12222382
$ws.Range("H21").Value = 7920
$ws.Range("J21").Value = 7920
$ws.Range("L21").Value = 7920
$ws.Range("N21").Value = -8266
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H24").Value = 8500
$ws.Range("J24").Value = 8500
$ws.Range("L24").Value = 8500
$ws.Range("N24").Value = -8846
$ws.Range("H29").Value = 5450
$ws.Range("J29").Value = 5450
$ws.Range("L29").Value = 5450
$ws.Range("N29").Value = -6030
$ws.Range("H30").Value = 7920
$ws.Range("J30").Value = 7920
$ws.Range("L30").Value = 7920
$ws.Range("N30").Value = -8130
$ws.Range("H122").Value = 3200
$ws.Range("I122").Value = 3760
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 11280
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -8830
$ws.Range("N122").Value = -12400
$ws.Range("H126").Value = 4004
$ws.Range("I126").Value = 2012
$ws.Range("K126").Value = 6036
$ws.Range("M126").Value = -3566
$ws.Range("H127").Value = 69326
$ws.Range("J127").Value = 69326
$ws.Range("L127").Value = 69326
$ws.Range("N127").Value = -79246
$ws.Range("H132").Value = 38468336
$ws.Range("I132").Value = 66675650
$ws.Range("J132").Value = 3820.4546
$ws.Range("K132").Value = 200026950
$ws.Range("L132").Value = 11461.3638
$ws.Range("M132").Value = -200024420
$ws.Range("N132").Value = -16521.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10501.25
$ws.Range("I7").Value = 10000
$ws.Range("K7").Value = 10000
$ws.Range("M7").Value = -9888
$ws.Range("H40").Value = 4113.875
$ws.Range("I40").Value = 4318.5
$ws.Range("J40").Value = 3500
$ws.Range("K40").Value = 4318.5
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = -4182.5
$ws.Range("N40").Value = -3772
$ws.Range("H126").Value = 10501.25
$ws.Range("I126").Value = 10000
$ws.Range("K126").Value = 30000
$ws.Range("M126").Value = -27530
$ws.Range("H132").Value = 2917.925
$ws.Range("I132").Value = 2165.0356
$ws.Range("J132").Value = 4674.6665
$ws.Range("K132").Value = 6495.1068
$ws.Range("L132").Value = 14023.9995
$ws.Range("M132").Value = -3965.1068
$ws.Range("N132").Value = -19083.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 900
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H126").Value = 1678.125
$ws.Range("I126").Value = 1656.6666
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 4969.9998
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2499.9998
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 5505392
$ws.Range("I132").Value = 2085.9524
$ws.Range("J132").Value = 26518016
$ws.Range("K132").Value = 6257.8572
$ws.Range("L132").Value = 79554048
$ws.Range("M132").Value = -3727.8572
$ws.Range("N132").Value = -79559108
$ws.Range("H136").Value = 2722.641
$ws.Range("I136").Value = 2593.4666
$ws.Range("J136").Value = 3153.2222
$ws.Range("K136").Value = 7780.399800000001
$ws.Range("L136").Value = 9459.6666
$ws.Range("M136").Value = -5230.399800000001
$ws.Range("N136").Value = -14559.6666
